# Update "想去人数" (want-to-go count) values in column F across sheets
# to reflect newly generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7119
$ws1.Range("F13").Value = 1442
$ws1.Range("F17").Value = 1147
$ws1.Range("F38").Value = 401
$ws1.Range("F39").Value = 61

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F16").Value = 1729
$ws2.Range("F26").Value = 621
$ws2.Range("F42").Value = 5

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value = 1438
$ws3.Range("F9").Value = 2261

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 7119
$ws4.Range("F15").Value = 1442
$ws4.Range("F18").Value = 1729
$ws4.Range("F21").Value = 1438
$ws4.Range("F22").Value = 2261
$ws4.Range("F25").Value = 1147
$ws4.Range("F31").Value = 621
$ws4.Range("F44").Value = 61
$ws4.Range("F48").Value = 5
